# Fixed concurrent SQLite bug, in accordance with
# https://github.com/jonathanburrows/lavalav/issues/22
#
# - Adds a new backlog item describing a fix for the SQLite concurrency bug
#   (new shared string, row 53 col G).
# - Marks that item (row 53, "Repository SQLite Bug Fixes") and its sibling
#   (row 54, "Assembly Loader Bug Fixes") half-complete (Points 0 -> 0.5).
# - Marks the "Change lvl.* folder" cleanup items (rows 41-49) and the
#   "Create Core Solution folder" item (row 51) as Completed.
# - Moves the sheet's scroll position / selection down towards the new work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backlog")
$ws.Activate()

# Rows 41-49: Status column (E) To Do -> Completed
$ws.Range("E41:E49").Value = "Completed"

# Row 51: Status column (E) In Progress -> Completed
$ws.Range("E51").Value = "Completed"

# Row 53 ("Repository SQLite Bug Fixes"): half-complete, with a description
# of the fix that was made.
$ws.Range("C53").Value = 0.5
$ws.Range("G53").Value = "In the repository, insides transactions; add a lock to the transactions, so that asynchronous SQLite operations don’t end transactions before complete."

# Row 54 ("Assembly Loader Bug Fixes"): half-complete too.
$ws.Range("C54").Value = 0.5

# Update the sheet's scroll position / active selection to reflect where
# work was happening.
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("G54").Select() | Out-Null
